$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "00061382"
$ws.Range("F1").Value = "Benjamin Munoz"

$ws.Range("B3").Value = "300006827"
$ws.Range("F3").Value = "Benjamin Munoz"

$ws.Range("B5").Value = "0880011926"
$ws.Range("F5").Value = "01/24"

$ws.Range("E10").Value = "customer said the craft shut down and fell 100 ft to the ground during the rth, he thought it was a product issue, not his responsibility, so asked for a data analysis"
